$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for row 2
$ws.Range("G2").Value = 22878
$ws.Range("H2").Value = 21072
$ws.Range("J2").Value = 165.04
$ws.Range("K2").Value = 222

# Update data values for row 3
$ws.Range("G3").Value = 22878
$ws.Range("H3").Value = 21072
$ws.Range("J3").Value = 165.04
$ws.Range("K3").Value = 222

# Update the selection shown in the sheet view
$ws.Range("K7").Select()
